# Apply updated crypto price/volume data (GitHub Actions refresh).
# Only the cells whose content actually changed are touched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '''62.616.48'
$ws.Cells.Item(2, 5).Value = '''  +3.23%  '
$ws.Cells.Item(3, 4).Value = '''2.446.10'
$ws.Cells.Item(3, 5).Value = '''  +2.00%  '
$ws.Cells.Item(4, 5).Value = '''  -0.08%  '
$ws.Cells.Item(5, 4).Value = '''578.70'
$ws.Cells.Item(6, 4).Value = '''145.46'
$ws.Cells.Item(6, 5).Value = '''  +3.11%  '
$ws.Cells.Item(7, 5).Value = '''  +0.08%  '
$ws.Cells.Item(8, 5).Value = '''  +0.62%  '
$ws.Cells.Item(9, 4).Value = '''2.444.39'
$ws.Cells.Item(9, 5).Value = '''  +1.71%  '
$ws.Cells.Item(10, 5).Value = '''  +2.39%  '
$ws.Cells.Item(11, 5).Value = '''  +1.06%  '
$ws.Cells.Item(12, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(12, 4).Value = '''0.0₅0168'
$ws.Cells.Item(12, 5).Value = '''  +504.52%  '
$ws.Cells.Item(13, 2).Value = 'Toncoin'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(13, 4).Value = '''5.23'
$ws.Cells.Item(13, 5).Value = '''  +1.45%  '
$ws.Cells.Item(14, 2).Value = 'Cardano'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Cells.Item(14, 4).Value = '''0.353'
$ws.Cells.Item(14, 5).Value = '''  +3.60%  '
$ws.Cells.Item(15, 2).Value = 'Avalanche'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(15, 4).Value = '''28.43'
$ws.Cells.Item(15, 5).Value = '''  +9.11%  '
$ws.Cells.Item(16, 2).Value = 'ShibaInu'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(16, 4).Value = '''0.0000178'
$ws.Cells.Item(16, 5).Value = '''  +5.86%  '
$ws.Cells.Item(17, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(17, 4).Value = '''2.889.10'
$ws.Cells.Item(17, 5).Value = '''  +1.87%  '
$ws.Cells.Item(18, 2).Value = 'WrappedBTC'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(18, 4).Value = '''62.585.45'
$ws.Cells.Item(18, 5).Value = '''  +3.24%  '
$ws.Cells.Item(19, 4).Value = '''2.437.08'
$ws.Cells.Item(19, 5).Value = '''  +1.05%  '
$ws.Cells.Item(20, 5).Value = '''  -1.25%  '
$ws.Cells.Item(21, 5).Value = '''  +2.64%  '
$ws.Cells.Item(22, 4).Value = '''326.05'
$ws.Cells.Item(22, 5).Value = '''  +0.81%  '
$ws.Cells.Item(24, 5).Value = '''  +11.32%  '
$ws.Cells.Item(25, 5).Value = '''  -0.03%  '
$ws.Cells.Item(26, 4).Value = '''65.37'
$ws.Cells.Item(26, 5).Value = '''  +0.61%  '
$ws.Cells.Item(27, 4).Value = '''646.04'
$ws.Cells.Item(27, 5).Value = '''  +14.82%  '
$ws.Cells.Item(28, 5).Value = '''  +14.91%  '
$ws.Cells.Item(29, 4).Value = '''8.55'
$ws.Cells.Item(29, 5).Value = '''  +6.31%  '
$ws.Cells.Item(30, 5).Value = '''  +4.84%  '
$ws.Cells.Item(31, 4).Value = '''2.562.05'
$ws.Cells.Item(31, 5).Value = '''  +1.96%  '
$ws.Cells.Item(32, 5).Value = '''  +1.41%  '
$ws.Cells.Item(33, 5).Value = '''  +6.81%  '
$ws.Cells.Item(34, 5).Value = '''  +3.09%  '
$ws.Cells.Item(35, 4).Value = '''0.139'
$ws.Cells.Item(35, 5).Value = '''  +6.23%  '
$ws.Cells.Item(36, 5).Value = '''  +2.27%  '
$ws.Cells.Item(37, 5).Value = '''  +0.11%  '
$ws.Cells.Item(38, 5).Value = '''  +3.41%  '
$ws.Cells.Item(39, 4).Value = '''154.04'
$ws.Cells.Item(39, 5).Value = '''  +1.22%  '
$ws.Cells.Item(40, 4).Value = '''5.46'
$ws.Cells.Item(40, 5).Value = '''  +6.54%  '
$ws.Cells.Item(41, 5).Value = '''  +0.82%  '
$ws.Cells.Item(42, 5).Value = '''  +1.83%  '
$ws.Cells.Item(43, 5).Value = '''  +8.47%  '
$ws.Cells.Item(44, 5).Value = '''  +4.92%  '
$ws.Cells.Item(45, 4).Value = '''42.58'
$ws.Cells.Item(45, 5).Value = '''  +2.17%  '
$ws.Cells.Item(46, 4).Value = '''0.999'
$ws.Cells.Item(46, 5).Value = '''  +0.02%  '
$ws.Cells.Item(48, 4).Value = '''144.33'
$ws.Cells.Item(48, 5).Value = '''  +2.13%  '
$ws.Cells.Item(49, 5).Value = '''  +1.57%  '
$ws.Cells.Item(50, 2).Value = 'ApolloX'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/WMDlfMJ1W+apollox-apx'
$ws.Cells.Item(50, 4).Value = '''0.207'
$ws.Cells.Item(50, 5).Value = '''  +431.35%  '
$ws.Cells.Item(51, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(51, 4).Value = '''20.53'
$ws.Cells.Item(51, 5).Value = '''  +6.90%  '
